$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, shifting existing rows 198-204 down to 199-205.
$ws.Rows.Item(198).EntireRow.Insert()

# Populate the newly inserted row 198 with the new weekly price record.
$ws.Range("A198").Value = 11
$ws.Range("B198").Value = "Vega Monumental Concepción"
$ws.Range("C198").Value = "Bíobío"
$ws.Range("D198").Value = 45041
$ws.Range("E198").Value = 8
$ws.Range("F198").Value = 100112043
$ws.Range("G198").Value = "Pepino ensalada"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 180
$ws.Range("K198").Value = 10000
$ws.Range("L198").Value = 11000
$ws.Range("M198").Value = 10556
$ws.Range("N198").Value = "$/caja 60 unidades"
$ws.Range("O198").Value = "Región de Arica y Parinacota"
$ws.Range("P198").Value = 176
$ws.Range("Q198").Value = 60
$ws.Range("R198").Value = "Hortaliza"
